$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the pre-edit values for columns C:F (runs, balls, fours, sixes),
# rows 2-12, before any cell gets overwritten.
$original = @{}
for ($r = 2; $r -le 12; $r++) {
    $original[$r] = @(
        $ws.Cells.Item($r, 3).Value2,
        $ws.Cells.Item($r, 4).Value2,
        $ws.Cells.Item($r, 5).Value2,
        $ws.Cells.Item($r, 6).Value2
    )
}

# The match rows got reordered (same player/team, stats permuted across
# innings). Mapping: destination row -> source row (pre-edit).
$mapping = @{
    2  = 5
    3  = 12
    4  = 11
    5  = 10
    6  = 3
    7  = 2
    8  = 4
    9  = 7
    10 = 9
    11 = 6
    12 = 8
}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $vals = $original[$oldRow]
    for ($col = 3; $col -le 6; $col++) {
        $text = [string]([int]$vals[$col - 3])
        $cell = $ws.Cells.Item($newRow, $col)
        # Leading apostrophe forces text storage (matching the sheet's
        # existing numbers-stored-as-text convention) while entering the
        # value the way a user typing it into Excel would.
        $cell.Value = "'" + $text
        # Drop the quote-prefix formatting flag COM applied for the above,
        # restoring the cell to the workbook's original unstyled look.
        $cell.ClearFormats()
    }
}

$wb.Save()
